$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to Text format first so dotted/decimal-looking
# values are preserved verbatim as strings instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.320.69"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "3.752.33"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "594.16"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Value = "169.42"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").Value = "3.748.22"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +4.16%  "
$ws.Range("D14").Value = "36.40"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "4.385.35"
$ws.Range("D16").Value = "3.760.57"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "18.53"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "67.346.31"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "7.18"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "10.48"
$ws.Range("E21").Value = "  -5.98%  "
$ws.Range("D22").Value = "466.08"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "0.717"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "83.79"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  -9.26%  "
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("D27").Value = "12.10"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("D30").Value = "2.89"
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("D31").Value = "3.906.77"
$ws.Range("D32").Value = "7.62"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "30.47"
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("D36").Value = "3.717.21"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "5.82"
$ws.Range("E41").Value = "  -2.11%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "0.310"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D45").Value = "8.70"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "45.81"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").Value = "395.62"
$ws.Range("E48").Value = "  -5.62%  "
$ws.Range("D49").Value = "0.000268"
$ws.Range("E49").Value = "  -8.79%  "
$ws.Range("D50").Value = "0.0352"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").Value = "38.92"
$ws.Range("E51").Value = "  +2.05%  "
